$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear cells that are removed entirely (columns J and O) for rows 2-5
$ws.Range("J2:J5").ClearContents()
$ws.Range("O2:O5").ClearContents()

# Row 2
$ws.Range("D2").Value = 5677
$ws.Range("E2").Value = 456
$ws.Range("F2").Value = 456
$ws.Range("G2").Value = -890
$ws.Range("H2").Value = -770
$ws.Range("I2").Value = -770
$ws.Range("K2").Value = 9849
$ws.Range("L2").Value = 8642
$ws.Range("M2").Value = 1207
$ws.Range("N2").Value = 1207
$ws.Range("P2").Value = 2226
$ws.Range("Q2").Value = 481
$ws.Range("R2").Value = -436
$ws.Range("S2").Value = 12
$ws.Range("T2").Value = 305
$ws.Range("U2").Value = 175
$ws.Range("V2").Value = 6584
$ws.Range("W2").Value = 8.029999999999999
$ws.Range("X2").Value = -13.56
$ws.Range("Y2").Value = -47.32
$ws.Range("Z2").Value = -7.54
$ws.Range("AA2").Value = 715.95
$ws.Range("AB2").Value = -7.67
$ws.Range("AC2").Value = -1730
$ws.Range("AD2").Value = -2.67
$ws.Range("AE2").Value = 2989
$ws.Range("AF2").Value = 1.54
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 44367832

# Row 3
$ws.Range("D3").Value = 6666
$ws.Range("E3").Value = 1250
$ws.Range("F3").Value = 1250
$ws.Range("G3").Value = 1158
$ws.Range("H3").Value = 1267
$ws.Range("I3").Value = 1267
$ws.Range("K3").Value = 9846
$ws.Range("L3").Value = 7306
$ws.Range("M3").Value = 2541
$ws.Range("N3").Value = 2541
$ws.Range("P3").Value = 2226
$ws.Range("Q3").Value = 1527
$ws.Range("R3").Value = -122
$ws.Range("S3").Value = -718
$ws.Range("T3").Value = 956
$ws.Range("U3").Value = 572
$ws.Range("V3").Value = 5257
$ws.Range("W3").Value = 18.75
$ws.Range("X3").Value = 19.01
$ws.Range("Y3").Value = 67.63
$ws.Range("Z3").Value = 12.87
$ws.Range("AA3").Value = 287.57
$ws.Range("AB3").Value = 46.15
$ws.Range("AC3").Value = 2847
$ws.Range("AD3").Value = 5.02
$ws.Range("AE3").Value = 5870
$ws.Range("AF3").Value = 2.44
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 44367832

# Row 4
$ws.Range("D4").Value = 7731
$ws.Range("E4").Value = 1724
$ws.Range("F4").Value = 1724
$ws.Range("G4").Value = 855
$ws.Range("H4").Value = 883
$ws.Range("I4").Value = 883
$ws.Range("K4").Value = 9868
$ws.Range("L4").Value = 6449
$ws.Range("M4").Value = 3419
$ws.Range("N4").Value = 3419
$ws.Range("P4").Value = 2226
$ws.Range("Q4").Value = 2483
$ws.Range("R4").Value = -1210
$ws.Range("S4").Value = -1185
$ws.Range("T4").Value = 923
$ws.Range("U4").Value = 1561
$ws.Range("V4").Value = 4223
$ws.Range("W4").Value = 22.3
$ws.Range("X4").Value = 11.43
$ws.Range("Y4").Value = 29.65
$ws.Range("Z4").Value = 8.960000000000001
$ws.Range("AA4").Value = 188.6
$ws.Range("AB4").Value = 87.08
$ws.Range("AC4").Value = 1985
$ws.Range("AD4").Value = 8.01
$ws.Range("AE4").Value = 7900
$ws.Range("AF4").Value = 2.01
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 44367832

# Row 5
$ws.Range("D5").Value = 6797
$ws.Range("E5").Value = 1432
$ws.Range("F5").Value = 1432
$ws.Range("G5").Value = 1050
$ws.Range("H5").Value = 1102
$ws.Range("I5").Value = 1102
$ws.Range("K5").Value = 9940
$ws.Range("L5").Value = 5382
$ws.Range("M5").Value = 4558
$ws.Range("N5").Value = 4558
$ws.Range("P5").Value = 2226
$ws.Range("Q5").Value = 1831
$ws.Range("R5").Value = -1019
$ws.Range("S5").Value = -858
$ws.Range("T5").Value = 759
$ws.Range("U5").Value = 1072
$ws.Range("V5").Value = 3556
$ws.Range("W5").Value = 21.07
$ws.Range("X5").Value = 16.21
$ws.Range("Y5").Value = 27.62
$ws.Range("Z5").Value = 11.13
$ws.Range("AA5").Value = 118.07
$ws.Range("AB5").Value = 137.66
$ws.Range("AC5").Value = 2475
$ws.Range("AD5").Value = 4.85
$ws.Range("AE5").Value = 10531
$ws.Range("AF5").Value = 1.14
$ws.Range("AG5").Value = 250
$ws.Range("AH5").Value = 2.08
$ws.Range("AI5").Value = 9.949999999999999
$ws.Range("AJ5").Value = 44367832

# Row 6
$ws.Range("D6").Value = 6693
$ws.Range("E6").Value = 1130
$ws.Range("F6").Value = 1130
$ws.Range("G6").Value = 961
$ws.Range("H6").Value = 868
$ws.Range("I6").Value = 868
$ws.Range("K6").Value = 10665
$ws.Range("L6").Value = 5081
$ws.Range("M6").Value = 5584
$ws.Range("N6").Value = 5584
$ws.Range("P6").Value = 2226
$ws.Range("Q6").Value = 1803
$ws.Range("R6").Value = -1173
$ws.Range("S6").Value = -509
$ws.Range("T6").Value = 1095
$ws.Range("U6").Value = 708
$ws.Range("V6").Value = 3050
$ws.Range("W6").Value = 16.88
$ws.Range("X6").Value = 12.97
$ws.Range("Y6").Value = 17.12
$ws.Range("Z6").Value = 8.43
$ws.Range("AA6").Value = 91
$ws.Range("AB6").Value = 186.48
$ws.Range("AC6").Value = 1951
$ws.Range("AD6").Value = 5.56
$ws.Range("AE6").Value = 12863
$ws.Range("AF6").Value = 0.84
$ws.Range("AG6").Value = 250
$ws.Range("AH6").Value = 2.3
$ws.Range("AI6").Value = 12.5
$ws.Range("AJ6").Value = 44398588

# Row 7
$ws.Range("D7").Value = 8087
$ws.Range("E7").Value = 1811
$ws.Range("G7").Value = 1556
$ws.Range("H7").Value = 1265
$ws.Range("I7").Value = 1137
$ws.Range("K7").Value = 11997
$ws.Range("L7").Value = 5256
$ws.Range("M7").Value = 6741
$ws.Range("N7").Value = 6741
$ws.Range("P7").Value = 2226
$ws.Range("Q7").Value = 1877
$ws.Range("R7").Value = -1097
$ws.Range("S7").Value = -250
$ws.Range("T7").Value = 967
$ws.Range("U7").Value = 1089
$ws.Range("W7").Value = 22.39
$ws.Range("X7").Value = 15.64
$ws.Range("Y7").Value = 18.45
$ws.Range("Z7").Value = 11.16
$ws.Range("AA7").Value = 77.97
$ws.Range("AC7").Value = 2554
$ws.Range("AD7").Value = 10.71
$ws.Range("AE7").Value = 15522
$ws.Range("AF7").Value = 1.76
$ws.Range("AG7").Value = 250
$ws.Range("AH7").Value = 0.91
$ws.Range("AI7").Value = 9.76

# Row 8
$ws.Range("D8").Value = 8585
$ws.Range("E8").Value = 1951
$ws.Range("G8").Value = 1679
$ws.Range("H8").Value = 1295
$ws.Range("I8").Value = 1364
$ws.Range("K8").Value = 13111
$ws.Range("L8").Value = 5184
$ws.Range("M8").Value = 7927
$ws.Range("N8").Value = 7927
$ws.Range("P8").Value = 2226
$ws.Range("Q8").Value = 1962
$ws.Range("R8").Value = -1031
$ws.Range("S8").Value = -219
$ws.Range("T8").Value = 945
$ws.Range("U8").Value = 1278
$ws.Range("W8").Value = 22.73
$ws.Range("X8").Value = 15.08
$ws.Range("Y8").Value = 18.6
$ws.Range("Z8").Value = 10.32
$ws.Range("AA8").Value = 65.40000000000001
$ws.Range("AC8").Value = 3064
$ws.Range("AD8").Value = 8.93
$ws.Range("AE8").Value = 18253
$ws.Range("AF8").Value = 1.5
$ws.Range("AG8").Value = 250
$ws.Range("AH8").Value = 0.91
$ws.Range("AI8").Value = 8.140000000000001

# Row 9
$ws.Range("D9").Value = 8766
$ws.Range("E9").Value = 1955
$ws.Range("G9").Value = 1764
$ws.Range("H9").Value = 1343
$ws.Range("I9").Value = 1343
$ws.Range("K9").Value = 14210
$ws.Range("L9").Value = 5048
$ws.Range("M9").Value = 9162
$ws.Range("N9").Value = 9162
$ws.Range("P9").Value = 2226
$ws.Range("Q9").Value = 2054
$ws.Range("R9").Value = -1062
$ws.Range("S9").Value = -191
$ws.Range("T9").Value = 992
$ws.Range("U9").Value = 1331
$ws.Range("W9").Value = 22.3
$ws.Range("X9").Value = 15.32
$ws.Range("Y9").Value = 15.72
$ws.Range("Z9").Value = 9.83
$ws.Range("AA9").Value = 55.1
$ws.Range("AC9").Value = 3017
$ws.Range("AD9").Value = 9.06
$ws.Range("AE9").Value = 21097
$ws.Range("AF9").Value = 1.3
$ws.Range("AG9").Value = 250
$ws.Range("AH9").Value = 0.91
$ws.Range("AI9").Value = 8.27
